$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.156.47"
$ws.Range("E2").Value = "  +3.70%  "

# Row 3
$ws.Range("D3").Value = "1.603.56"
$ws.Range("E3").Value = "  +3.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "212.67"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "

# Row 6
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
$ws.Range("E7").Value = "  +2.29%  "

# Row 8
$ws.Range("E8").Value = "  +2.82%  "

# Row 9
$ws.Range("E9").Value = "  +1.39%  "

# Row 10
$ws.Range("E10").Value = "  +1.24%  "

# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0817"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +4.88%  "

# Row 12
$ws.Range("D12").Value = "1.827.12"
$ws.Range("E12").Value = "  +3.55%  "

# Row 13
$ws.Range("D13").Value = "1.606.27"
$ws.Range("E13").Value = "  +3.76%  "

# Row 14
$ws.Range("E14").Value = "  +0.59%  "

# Row 15
$ws.Range("E15").Value = "  +1.52%  "

# Row 16
$ws.Range("D16").Value = "26.146.60"
$ws.Range("E16").Value = "  +3.79%  "

# Row 17
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "60.47"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +3.31%  "

# Row 18
$ws.Range("E18").Value = "  +2.21%  "

# Row 19
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "204.35"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +10.09%  "

# Row 21
$ws.Range("E21").Value = "  +3.46%  "

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "9.30"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23
$ws.Range("E23").Value = "  +2.82%  "

# Row 24
$ws.Range("E24").Value = "  +11.08%  "

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "141.51"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +1.45%  "

# Row 26
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("E27").Value = "  -3.15%  "

# Row 28
$ws.Range("E28").Value = "  +2.72%  "

# Row 29
$ws.Range("E29").Value = "  +0.68%  "

# Row 30
$ws.Range("E30").Value = "  +1.57%  "

# Row 31
$ws.Range("E31").Value = "  +2.03%  "

# Row 32
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.11"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +3.39%  "

# Row 33
$ws.Range("E33").Value = "  +0.58%  "

# Row 34
$ws.Range("E34").Value = "  +1.64%  "

# Row 35
$ws.Range("E35").Value = "  +1.93%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.118.26"
$ws.Range("E36").Value = "  +3.00%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.0164"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +9.74%  "

# Row 38
$ws.Range("E38").Value = "  -0.22%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.783"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +3.27%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "2.30"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +2.55%  "

# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.491"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -0.59%  "

# Row 42
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.782"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "

# Row 43
$ws.Range("D43").Value = "1.739.42"
$ws.Range("E43").Value = "  +3.55%  "

# Row 44
$ws.Range("E44").Value = "  +1.88%  "

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "92.96"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +0.40%  "

# Row 46
$ws.Range("E46").Value = "  +4.07%  "

# Row 47
$ws.Range("E47").Value = "  +2.37%  "

# Row 48
$ws.Range("E48").Value = "  +0.57%  "

# Row 49
$ws.Range("E49").Value = "  +1.15%  "

# Row 50
$ws.Range("E50").Value = "  -0.02%  "

# Row 51
$ws.Range("D51").Value = "0.0₇0925"
$ws.Range("E51").Value = "  -15.93%  "
